# "Logged Week 15 and simulated Week 16"
# Update the Row ("R") totals for both the OFF and DEF sheets with the
# latest weekly cumulative stats.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row labeled "R" (row 3) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 200
$wsOff.Range("C3").Value = 147
$wsOff.Range("D3").Value = 51
$wsOff.Range("E3").Value = 23

# --- DEF sheet: row labeled "R" (row 3) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 287
$wsDef.Range("C3").Value = 211
$wsDef.Range("D3").Value = 57
$wsDef.Range("E3").Value = 32
$wsDef.Range("F3").Value = 4
